$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows before row 80, pushing existing rows 80-148 down to 88-156.
$ws.Range("A80:A87").EntireRow.Insert()

# Fixed (constant) values shared by every data row in this block.
$fA = 2
$fB = 'Comercializadora del Agro de Limarí'
$fC = 'Coquimbo'
$fE = 4
$fF = 'Fruta'
$fG = 100103
$fH = 'Frutos de hueso (carozo)'
$fI = 100103001
$fJ = 'Cereza'

# New weekly data rows (80-87): Date, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, Precio$/Kg, Kg/unidad
$newRows = @(
    @(44574,'Lapins','Especial',400,10500,11000,10750,'$/bandeja 10 kilos','Provincia de Curicó',1075,10),
    @(44574,'Lapins','Primera',300,8500,9000,8750,'$/bandeja 10 kilos','Provincia de Curicó',875,10),
    @(44574,'Lapins','Segunda',300,6500,7000,6750,'$/bandeja 10 kilos','Provincia de Curicó',675,10),
    @(44574,'Santina','Especial',360,10500,11000,10750,'$/bandeja 10 kilos','Provincia de Curicó',1075,10),
    @(44574,'Santina','Primera',400,8500,9000,8750,'$/bandeja 10 kilos','Provincia de Curicó',875,10),
    @(44574,'Santina','Segunda',300,6500,7000,6750,'$/bandeja 10 kilos','Provincia de Curicó',675,10),
    @(44574,'Sweet Heart','Especial',360,10500,11000,10750,'$/bandeja 10 kilos','Provincia de Curicó',1075,10),
    @(44574,'Sweet Heart','Primera',480,8500,9000,8708,'$/bandeja 10 kilos','Provincia de Curicó',871,10)
)

$startRow = 80
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $d = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $fA
    $ws.Cells.Item($r, 2).Value = $fB
    $ws.Cells.Item($r, 3).Value = $fC
    $ws.Cells.Item($r, 4).Value = $d[0]
    $ws.Cells.Item($r, 5).Value = $fE
    $ws.Cells.Item($r, 6).Value = $fF
    $ws.Cells.Item($r, 7).Value = $fG
    $ws.Cells.Item($r, 8).Value = $fH
    $ws.Cells.Item($r, 9).Value = $fI
    $ws.Cells.Item($r, 10).Value = $fJ
    $ws.Cells.Item($r, 11).Value = $d[1]
    $ws.Cells.Item($r, 12).Value = $d[2]
    $ws.Cells.Item($r, 13).Value = $d[3]
    $ws.Cells.Item($r, 14).Value = $d[4]
    $ws.Cells.Item($r, 15).Value = $d[5]
    $ws.Cells.Item($r, 16).Value = $d[6]
    $ws.Cells.Item($r, 17).Value = $d[7]
    $ws.Cells.Item($r, 18).Value = $d[8]
    $ws.Cells.Item($r, 19).Value = $d[9]
    $ws.Cells.Item($r, 20).Value = $d[10]
}
